$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 3003
    $ws.Range("F8").Value = 2059
    $ws.Range("F11").Value = 905
}
